$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(95, 8).Value = 32997.332
$ws.Cells.Item(95, 10).Value = 32997.332
$ws.Cells.Item(95, 12).Value = 32997.332
$ws.Cells.Item(95, 14).Value = -38489.332
$ws.Cells.Item(98, 8).Value = 31455.076
$ws.Cells.Item(98, 9).Value = 1116.5625
$ws.Cells.Item(98, 10).Value = 79996.7
$ws.Cells.Item(98, 11).Value = 1116.5625
$ws.Cells.Item(98, 12).Value = 79996.7
$ws.Cells.Item(98, 13).Value = 381.4375
$ws.Cells.Item(98, 14).Value = -82992.7
$ws.Cells.Item(105, 8).Value = 38661.332
$ws.Cells.Item(105, 10).Value = 38661.332
$ws.Cells.Item(105, 12).Value = 38661.332
$ws.Cells.Item(105, 14).Value = -45649.332
$ws.Cells.Item(122, 8).Value = 31455.076
$ws.Cells.Item(122, 9).Value = 1116.5625
$ws.Cells.Item(122, 10).Value = 79996.7
$ws.Cells.Item(122, 11).Value = 3349.6875
$ws.Cells.Item(122, 12).Value = 239990.1
$ws.Cells.Item(122, 13).Value = -899.6875
$ws.Cells.Item(122, 14).Value = -244890.1
$ws.Cells.Item(129, 8).Value = 5260.857
$ws.Cells.Item(129, 10).Value = 5695.5
$ws.Cells.Item(129, 12).Value = 17086.5
$ws.Cells.Item(129, 14).Value = -27086.5
$ws.Cells.Item(138, 8).Value = 2867.35
$ws.Cells.Item(138, 9).Value = 2439.3333
$ws.Cells.Item(138, 10).Value = 3152.6943
$ws.Cells.Item(138, 11).Value = 7317.999899999999
$ws.Cells.Item(138, 12).Value = 9458.082900000001
$ws.Cells.Item(138, 13).Value = -2177.999899999999
$ws.Cells.Item(138, 14).Value = -19738.0829

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(31, 8).Value = 10760
$ws.Cells.Item(31, 9).Value = 3368
$ws.Cells.Item(31, 11).Value = 3368
$ws.Cells.Item(31, 13).Value = -3074
$ws.Cells.Item(32, 8).Value = 29272.896
$ws.Cells.Item(32, 9).Value = 30905.416
$ws.Cells.Item(32, 10).Value = 15279.857
$ws.Cells.Item(32, 11).Value = 30905.416
$ws.Cells.Item(32, 12).Value = 15279.857
$ws.Cells.Item(32, 13).Value = -30618.416
$ws.Cells.Item(32, 14).Value = -15853.857
$ws.Cells.Item(61, 8).Value = 3126.16
$ws.Cells.Item(61, 9).Value = 3075.2273
$ws.Cells.Item(61, 11).Value = 3075.2273
$ws.Cells.Item(61, 13).Value = -2863.2273
$ws.Cells.Item(95, 8).Value = 27604
$ws.Cells.Item(95, 10).Value = 27604
$ws.Cells.Item(95, 12).Value = 27604
$ws.Cells.Item(95, 14).Value = -33096
$ws.Cells.Item(101, 8).Value = 49400.4
$ws.Cells.Item(101, 10).Value = 49400.4
$ws.Cells.Item(101, 12).Value = 49400.4
$ws.Cells.Item(101, 14).Value = -55890.4
$ws.Cells.Item(103, 8).Value = 35695.332
$ws.Cells.Item(103, 10).Value = 35695.332
$ws.Cells.Item(103, 12).Value = 35695.332
$ws.Cells.Item(103, 14).Value = -38039.332
$ws.Cells.Item(104, 8).Value = 30644.8
$ws.Cells.Item(104, 10).Value = 30644.8
$ws.Cells.Item(104, 12).Value = 30644.8
$ws.Cells.Item(104, 14).Value = -37632.8
$ws.Cells.Item(105, 8).Value = 48996
$ws.Cells.Item(105, 10).Value = 48996
$ws.Cells.Item(105, 12).Value = 48996
$ws.Cells.Item(105, 14).Value = -55984
$ws.Cells.Item(106, 8).Value = 48368.668
$ws.Cells.Item(106, 10).Value = 48368.668
$ws.Cells.Item(106, 12).Value = 48368.668
$ws.Cells.Item(106, 14).Value = -50892.668
$ws.Cells.Item(132, 8).Value = 10418721
$ws.Cells.Item(132, 9).Value = 12501779
$ws.Cells.Item(132, 10).Value = 3429.5
$ws.Cells.Item(132, 11).Value = 37505337
$ws.Cells.Item(132, 12).Value = 10288.5
$ws.Cells.Item(132, 13).Value = -37502807
$ws.Cells.Item(132, 14).Value = -15348.5
$ws.Cells.Item(136, 8).Value = 3126.16
$ws.Cells.Item(136, 9).Value = 3075.2273
$ws.Cells.Item(136, 11).Value = 9225.6819
$ws.Cells.Item(136, 13).Value = -6675.6819

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(95, 8).Value = 44727.332
$ws.Cells.Item(95, 10).Value = 44727.332
$ws.Cells.Item(95, 12).Value = 44727.332
$ws.Cells.Item(95, 14).Value = -50219.332
$ws.Cells.Item(124, 8).Value = 44647.668
$ws.Cells.Item(124, 10).Value = 44647.668
$ws.Cells.Item(124, 12).Value = 44647.668
$ws.Cells.Item(124, 14).Value = -54467.668
$ws.Cells.Item(125, 8).Value = 0
$ws.Cells.Item(125, 10).Value = 0
$ws.Cells.Item(125, 12).Value = 0
$ws.Cells.Item(125, 14).ClearContents()
$ws.Cells.Item(126, 8).Value = 50780
$ws.Cells.Item(126, 10).Value = 50780
$ws.Cells.Item(126, 12).Value = 50780
$ws.Cells.Item(126, 14).Value = -60660
$ws.Cells.Item(129, 8).Value = 33333
$ws.Cells.Item(129, 10).Value = 33333
$ws.Cells.Item(129, 12).Value = 33333
$ws.Cells.Item(129, 14).Value = -43333
$ws.Cells.Item(130, 8).Value = 50413
$ws.Cells.Item(130, 10).Value = 50413
$ws.Cells.Item(130, 12).Value = 50413
$ws.Cells.Item(130, 14).Value = -60453
$ws.Cells.Item(134, 8).Value = 3724.3157
$ws.Cells.Item(134, 9).Value = 3836.5
$ws.Cells.Item(134, 10).Value = 3642.7273
$ws.Cells.Item(134, 11).Value = 11509.5
$ws.Cells.Item(134, 12).Value = 10928.1819
$ws.Cells.Item(134, 13).Value = -8974.5
$ws.Cells.Item(134, 14).Value = -15998.1819

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(28, 8).Value = 35158.555
$ws.Cells.Item(28, 10).Value = 35158.555
$ws.Cells.Item(28, 12).Value = 35158.555
$ws.Cells.Item(28, 14).Value = -35648.555
$ws.Cells.Item(31, 8).Value = 8630.92
$ws.Cells.Item(31, 9).Value = 4432.625
$ws.Cells.Item(31, 10).Value = 10606.588
$ws.Cells.Item(31, 11).Value = 4432.625
$ws.Cells.Item(31, 12).Value = 10606.588
$ws.Cells.Item(31, 13).Value = -4137.625
$ws.Cells.Item(31, 14).Value = -11196.588
$ws.Cells.Item(34, 8).Value = 8630.92
$ws.Cells.Item(34, 9).Value = 4432.625
$ws.Cells.Item(34, 10).Value = 10606.588
$ws.Cells.Item(34, 11).Value = 4432.625
$ws.Cells.Item(34, 12).Value = 10606.588
$ws.Cells.Item(34, 13).Value = -4230.625
$ws.Cells.Item(34, 14).Value = -11010.588
$ws.Cells.Item(43, 8).Value = 28219
$ws.Cells.Item(43, 10).Value = 28219
$ws.Cells.Item(43, 12).Value = 28219
$ws.Cells.Item(43, 14).Value = -28587
$ws.Cells.Item(92, 8).Value = 31867
$ws.Cells.Item(92, 10).Value = 31867
$ws.Cells.Item(92, 12).Value = 31867
$ws.Cells.Item(92, 14).Value = -36859
$ws.Cells.Item(96, 8).Value = 71687
$ws.Cells.Item(96, 10).Value = 71687
$ws.Cells.Item(96, 12).Value = 71687
$ws.Cells.Item(96, 14).Value = -77179
$ws.Cells.Item(101, 8).Value = 28219
$ws.Cells.Item(101, 10).Value = 28219
$ws.Cells.Item(101, 12).Value = 28219
$ws.Cells.Item(101, 14).Value = -34709
$ws.Cells.Item(106, 8).Value = 47992
$ws.Cells.Item(106, 10).Value = 47992
$ws.Cells.Item(106, 12).Value = 47992
$ws.Cells.Item(106, 14).Value = -50516
$ws.Cells.Item(122, 8).Value = 71516720
$ws.Cells.Item(122, 9).Value = 91019150
$ws.Cells.Item(122, 11).Value = 273057450
$ws.Cells.Item(122, 13).Value = -273055000
$ws.Cells.Item(127, 8).Value = 0
$ws.Cells.Item(127, 10).Value = 0
$ws.Cells.Item(127, 12).Value = 0
$ws.Cells.Item(127, 14).ClearContents()
$ws.Cells.Item(132, 8).Value = 42294.258
$ws.Cells.Item(132, 9).Value = 1600.3478
$ws.Cells.Item(132, 10).Value = 120290.914
$ws.Cells.Item(132, 11).Value = 4801.0434
$ws.Cells.Item(132, 12).Value = 360872.742
$ws.Cells.Item(132, 13).Value = -2271.0434
$ws.Cells.Item(132, 14).Value = -365932.742

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(98, 8).Value = 801
$ws.Cells.Item(98, 9).Value = 750
$ws.Cells.Item(98, 10).Value = 852
$ws.Cells.Item(98, 11).Value = 2250
$ws.Cells.Item(98, 12).Value = 2556
$ws.Cells.Item(98, 13).Value = -752
$ws.Cells.Item(98, 14).Value = -5552
$ws.Cells.Item(131, 8).Value = 2280.4404
$ws.Cells.Item(131, 9).Value = 99999
$ws.Cells.Item(131, 10).Value = 1103.1084
$ws.Cells.Item(131, 11).Value = 299997
$ws.Cells.Item(131, 12).Value = 3309.3252
$ws.Cells.Item(131, 13).Value = -294957
$ws.Cells.Item(131, 14).Value = -13389.3252

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(132, 8).Value = 3169.4443
$ws.Cells.Item(132, 9).Value = 2541.4
$ws.Cells.Item(132, 10).Value = 3954.5
$ws.Cells.Item(132, 11).Value = 7624.200000000001
$ws.Cells.Item(132, 12).Value = 11863.5
$ws.Cells.Item(132, 13).Value = -5094.200000000001
$ws.Cells.Item(132, 14).Value = -16923.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(105, 8).Value = 0
$ws.Cells.Item(105, 10).Value = 0
$ws.Cells.Item(105, 12).Value = 0
$ws.Cells.Item(105, 14).ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(95, 8).Value = 34086
$ws.Cells.Item(95, 10).Value = 34086
$ws.Cells.Item(95, 12).Value = 34086
$ws.Cells.Item(95, 14).Value = -39578
$ws.Cells.Item(103, 8).Value = 41641.43
$ws.Cells.Item(103, 10).Value = 41641.43
$ws.Cells.Item(103, 12).Value = 41641.43
$ws.Cells.Item(103, 14).Value = -43985.43
$ws.Cells.Item(104, 8).Value = 39996
$ws.Cells.Item(104, 10).Value = 39996
$ws.Cells.Item(104, 12).Value = 39996
$ws.Cells.Item(104, 14).Value = -46984
$ws.Cells.Item(132, 8).Value = 1316.9574
$ws.Cells.Item(132, 9).Value = 1039.95
$ws.Cells.Item(132, 10).Value = 2899.8572
$ws.Cells.Item(132, 11).Value = 3119.85
$ws.Cells.Item(132, 12).Value = 8699.571599999999
$ws.Cells.Item(132, 13).Value = -589.8500000000004
$ws.Cells.Item(132, 14).Value = -13759.5716
